$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The original sheet has 9 columns:
#   A = running index (unused helper column)              -> dropped
#   B = "config" full repr of the model pipeline (wraptext)-> dropped
#   C = mse_oos
#   D = mse_validated
#   E = msfe_adj
#   F = name  (short model name: const, ols, pca, ...)
#   G = period
#   H = r2_oos
#   I = start_idx
#
# The new layout keeps only the 7 meaningful columns, reordered as:
#   A = name
#   B = mse_oos
#   C = mse_validated
#   D = r2_oos
#   E = msfe_adj
#   F = period
#   G = start_idx
# ---------------------------------------------------------------------------

$rowCount = $ws.UsedRange.Rows.Count

# Capture every value we still need before touching the sheet structure.
$mse_oos = @()
$mse_validated = @()
$msfe_adj = @()
$name = @()
$period = @()
$r2_oos = @()
$start_idx = @()

for ($r = 1; $r -le $rowCount; $r++) {
    $mse_oos += ,$ws.Cells.Item($r, 3).Value2
    $mse_validated += ,$ws.Cells.Item($r, 4).Value2
    $msfe_adj += ,$ws.Cells.Item($r, 5).Value2
    $name += ,$ws.Cells.Item($r, 6).Value2
    $period += ,$ws.Cells.Item($r, 7).Value2
    $r2_oos += ,$ws.Cells.Item($r, 8).Value2
    $start_idx += ,$ws.Cells.Item($r, 9).Value2
}

# Wipe all cell content/formatting, then drop the now unused index + config
# columns (this also clears the custom column-width set on the old column B).
$ws.Cells.Clear()
$ws.Columns("A:B").Delete()

# Write the header row back in the new order.
$ws.Cells.Item(1, 1).Value = $name[0]
$ws.Cells.Item(1, 2).Value = $mse_oos[0]
$ws.Cells.Item(1, 3).Value = $mse_validated[0]
$ws.Cells.Item(1, 4).Value = $r2_oos[0]
$ws.Cells.Item(1, 5).Value = $msfe_adj[0]
$ws.Cells.Item(1, 6).Value = $period[0]
$ws.Cells.Item(1, 7).Value = $start_idx[0]

# Write the data rows back in the new order.
for ($i = 1; $i -lt $rowCount; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $name[$i]
    $ws.Cells.Item($r, 2).Value2 = $mse_oos[$i]
    $ws.Cells.Item($r, 3).Value2 = $mse_validated[$i]
    $ws.Cells.Item($r, 4).Value2 = $r2_oos[$i]
    $ws.Cells.Item($r, 5).Value2 = $msfe_adj[$i]
    $ws.Cells.Item($r, 6).Value2 = $period[$i]
    $ws.Cells.Item($r, 7).Value2 = $start_idx[$i]
}

$ws.Application.GoTo($ws.Range("K15"))
